$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("G1").Value = "answer (optional column, for checking only)"
$ws.Range("H1").Value = "length (optional column, for checking only)"

# --- Answer / length data for each clue row (rows 2-11) ---
$answers = @("tab", "irl", "paris", "broth", "jazzy", "tiara", "arroz", "blitz", "pbj", "shy")
$lengths = @(3, 3, 5, 5, 5, 5, 5, 5, 3, 3)

for ($i = 0; $i -lt $answers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $answers[$i]
    $ws.Cells.Item($row, 8).Value = $lengths[$i]
}

# --- Column widths (approximate target pixel-grid widths as closely as the host allows) ---
$ws.Columns.Item(6).ColumnWidth = 29.333333333333332
$ws.Columns.Item(7).ColumnWidth = 23.5
$ws.Columns.Item(8).ColumnWidth = 35.333333333333336

# --- View state: zoom + active selection ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("H1").Select() | Out-Null
